# "Add 710 to electives list"
#
# Inserts a new elective row (EPIB 710 Advanced Causal Inference, 3 credits,
# Biostatistics category) just above the old row 30 ("PPHS 501 Population
# Health and Epidemiology"), which pushes it and everything below down by
# one row. Also refreshes the "Last edited" footer line with the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30 (existing row 30 onward shift down to 31+).
$ws.Rows.Item(30).Insert()

# Populate the new row with the new course.
$ws.Cells.Item(30, 1).Value = "EPIB 710 Advanced Causal Inference"
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(30, 3).Value = "Biostatistics"

# The footer "Last edited" note is now on the last row of the sheet; update
# it to reflect the new edit date/author.
$lastRow = $ws.UsedRange.Rows.Count
$ws.Cells.Item($lastRow, 1).Value = "Last edited: 2022-06-07 by Sam Harper"

# Reflect the view state that was left behind in the source file (purely
# cosmetic - scrolled down a bit with A52 selected/active).
$ws.Range("A52").Select()
$excel.ActiveWindow.ScrollRow = 18
